$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new "category" column between property_category (H) and date (old I)
# this shifts old I/J/K (date/legislator_name/legislator_id) one column right -> J/K/L
$ws.Columns.Item(9).Insert()

$ws.Range("I1").Value = "category"
$ws.Range("I2").Value = "normal"
$ws.Range("I3").Value = "normal"

# Append two new trailing columns: source_file, index
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"
$ws.Range("M2").Value = "tmped871"
$ws.Range("N2").Value = 52
$ws.Range("M3").Value = "tmped871"
$ws.Range("N3").Value = 53

# Match header styling (bold + border + centered) on the new trailing header cells
$ws.Range("B1").Copy()
$ws.Range("M1:N1").PasteSpecial(-4122)

# Match data-row styling on the new trailing data cells
$ws.Range("B2").Copy()
$ws.Range("M2:N2").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("M3:N3").PasteSpecial(-4122)

$excel.CutCopyMode = $false
